# Auto-generated Excel COM-interop script applying the eterno_data.xlsx diff
# Adds new rows to Users, Products, POS_Sales, Customer_Orders sheets,
# and updates two existing Products cells (row 2 and row 3).

$wb = $excel.ActiveWorkbook

# ---- Sheet: Users ----
$ws = $wb.Worksheets.Item('Users')
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'kaizen'
$ws.Cells.Item(4, 3).Value = 'boarratjabol@gmail.com'
$ws.Cells.Item(4, 4).Value = 'customer'
$ws.Cells.Item(4, 5).Value = '2025-11-09 13:59:58'

# ---- Sheet: Products ----
$ws = $wb.Worksheets.Item('Products')
$ws.Cells.Item(2, 2).Value = 'Eterno Void'
$ws.Cells.Item(2, 4).Value = 599
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 7).Value = 'https://i.ibb.co/sd7crdHV/468614642-122110597646602772-5741305816202520031-n.jpg'
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 7).Value = 'https://i.ibb.co/23s110xr/504703845-122141015084602772-3956461215141399652-n.jpg'
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'Eterno Outlaw'
$ws.Cells.Item(4, 4).Value = 1199
$ws.Cells.Item(4, 5).Value = 20
$ws.Cells.Item(4, 6).Value = 'Jackets'
$ws.Cells.Item(4, 7).Value = 'https://i.ibb.co/4RvHGQnL/504812987-122140885244602772-6478208165150186098-n.jpg'
$ws.Cells.Item(4, 8).Value = '2025-11-09 13:17:33'
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'Eterno Saint'
$ws.Cells.Item(5, 4).Value = 1399
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 'Jackets'
$ws.Cells.Item(5, 7).Value = 'https://i.ibb.co/fz1y42b0/545423967-122152407638602772-4284048291396969698-n.jpg'
$ws.Cells.Item(5, 8).Value = '2025-11-09 13:17:53'
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'Eterno Pulse'
$ws.Cells.Item(6, 4).Value = 499
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 'Shirts'
$ws.Cells.Item(6, 7).Value = 'https://i.ibb.co/GQVDzrxf/Screenshot-2025-11-09-212227.png'
$ws.Cells.Item(6, 8).Value = '2025-11-09 13:33:38'
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'Eterno Sanctum'
$ws.Cells.Item(7, 4).Value = 499
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 'Shirts'
$ws.Cells.Item(7, 7).Value = 'https://i.ibb.co/BHDjYFJT/504252726-122141015126602772-536880939443070759-n.jpg'
$ws.Cells.Item(7, 8).Value = '2025-11-09 13:34:32'
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'Eterno Eclipse'
$ws.Cells.Item(8, 4).Value = 4999
$ws.Cells.Item(8, 5).Value = 10
$ws.Cells.Item(8, 6).Value = 'Outerwear'
$ws.Cells.Item(8, 7).Value = 'https://i.ibb.co/27PxvGqV/image.jpg'
$ws.Cells.Item(8, 8).Value = '2025-11-09 13:49:55'
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Eterno Abyss'
$ws.Cells.Item(9, 4).Value = 4599
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 'Outerwear'
$ws.Cells.Item(9, 7).Value = 'https://i.ibb.co/356xp46C/Screenshot-2025-11-09-214504.png'
$ws.Cells.Item(9, 8).Value = '2025-11-09 13:50:42'
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'Eterno Drift'
$ws.Cells.Item(10, 4).Value = 3999
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(10, 6).Value = 'Outerwear'
$ws.Cells.Item(10, 7).Value = 'https://i.ibb.co/kV3YC222/Screenshot-2025-11-09-214421.png'
$ws.Cells.Item(10, 8).Value = '2025-11-09 13:51:23'
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Eterno Ethos'
$ws.Cells.Item(11, 4).Value = 1999
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 'Jackets'
$ws.Cells.Item(11, 7).Value = 'https://i.ibb.co/HcTcqJD/Screenshot-2025-11-09-214534.png'
$ws.Cells.Item(11, 8).Value = '2025-11-09 13:52:04'
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Bonest Gatti'
$ws.Cells.Item(12, 4).Value = 25799
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 'Accessories'
$ws.Cells.Item(12, 7).Value = 'https://i.ibb.co/jvqykhMB/Screenshot-2025-11-09-213836.png'
$ws.Cells.Item(12, 8).Value = '2025-11-09 13:53:39'
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'Bonest Gatti'
$ws.Cells.Item(13, 4).Value = 15999
$ws.Cells.Item(13, 5).Value = 5
$ws.Cells.Item(13, 6).Value = 'Accessories'
$ws.Cells.Item(13, 7).Value = 'https://i.ibb.co/DPXp1GYv/Screenshot-2025-11-09-213909.png'
$ws.Cells.Item(13, 8).Value = '2025-11-09 13:54:20'
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'Bonest Gatti'
$ws.Cells.Item(14, 4).Value = 26789
$ws.Cells.Item(14, 5).Value = 5
$ws.Cells.Item(14, 6).Value = 'Accessories'
$ws.Cells.Item(14, 7).Value = 'https://i.ibb.co/GvXy6tG1/Screenshot-2025-11-09-214106.png'
$ws.Cells.Item(14, 8).Value = '2025-11-09 13:55:09'
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'Bonest Gatti'
$ws.Cells.Item(15, 4).Value = 24509
$ws.Cells.Item(15, 5).Value = 5
$ws.Cells.Item(15, 6).Value = 'Accessories'
$ws.Cells.Item(15, 7).Value = 'https://i.ibb.co/dsVq1kRp/Screenshot-2025-11-09-214036.png'
$ws.Cells.Item(15, 8).Value = '2025-11-09 13:55:51'

# ---- Sheet: POS_Sales ----
$ws = $wb.Worksheets.Item('POS_Sales')
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 1998
$ws.Cells.Item(4, 4).Value = 'cash'
$ws.Cells.Item(4, 5).Value = 'voucher'
$ws.Cells.Item(4, 6).Value = 100
$ws.Cells.Item(4, 7).Value = '[{"product_id": 2, "name": "Eterno Grace", "price": 799, "quantity": 2, "stock": 19}, {"product_id": 1, "name": "Eterno Shirt", "price": 500, "quantity": 1, "stock": 5}]'
$ws.Cells.Item(4, 8).Value = '2025-11-08 07:25:20'
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = 20639.2
$ws.Cells.Item(5, 4).Value = 'cash'
$ws.Cells.Item(5, 5).Value = 'senior'
$ws.Cells.Item(5, 6).Value = 5159.8
$ws.Cells.Item(5, 7).Value = '[{"product_id": 11, "name": "Bonest Gatti", "price": 25799, "quantity": 1, "stock": 4}]'
$ws.Cells.Item(5, 8).Value = '2025-11-09 14:18:29'

# ---- Sheet: Customer_Orders ----
$ws = $wb.Worksheets.Item('Customer_Orders')
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = 'test'
$ws.Cells.Item(5, 4).Value = 'test@gmail.com'
$ws.Cells.Item(5, 5).Value = 'N/A'
$ws.Cells.Item(5, 6).Value = 1299
$ws.Cells.Item(5, 7).Value = 55
$ws.Cells.Item(5, 8).Value = 1354
$ws.Cells.Item(5, 9).Value = 'credit_card'
$ws.Cells.Item(5, 10).Value = 'completed'
$ws.Cells.Item(5, 11).Value = '[{"product_id": 2, "product_name": "Eterno Grace", "quantity": 1, "price": 799.0}, {"product_id": 1, "product_name": "Eterno Shirt", "quantity": 1, "price": 500.0}]'
$ws.Cells.Item(5, 12).Value = '2025-11-09 12:09:55'
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = 'test'
$ws.Cells.Item(6, 4).Value = 'test@gmail.com'
$ws.Cells.Item(6, 5).Value = 'N/A'
$ws.Cells.Item(6, 6).Value = 500
$ws.Cells.Item(6, 7).Value = 167
$ws.Cells.Item(6, 8).Value = 567
$ws.Cells.Item(6, 9).Value = 'paypal'
$ws.Cells.Item(6, 10).Value = 'pending'
$ws.Cells.Item(6, 11).Value = '[{"product_id": 1, "product_name": "Eterno Shirt", "quantity": 1, "price": 500.0}]'
$ws.Cells.Item(6, 12).Value = '2025-11-09 12:12:14'
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Value = 'test'
$ws.Cells.Item(7, 4).Value = 'test@gmail.com'
$ws.Cells.Item(7, 5).Value = 'N/A'
$ws.Cells.Item(7, 6).Value = 799
$ws.Cells.Item(7, 7).Value = 61
$ws.Cells.Item(7, 8).Value = 760
$ws.Cells.Item(7, 9).Value = 'gcash'
$ws.Cells.Item(7, 10).Value = 'pending'
$ws.Cells.Item(7, 11).Value = '[{"product_id": 2, "product_name": "Eterno Grace", "quantity": 1, "price": 799.0}]'
$ws.Cells.Item(7, 12).Value = '2025-11-09 12:12:59'
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 2
$ws.Cells.Item(8, 3).Value = 'test'
$ws.Cells.Item(8, 4).Value = 'test@gmail.com'
$ws.Cells.Item(8, 5).Value = 'egg street egg city egg municipal'
$ws.Cells.Item(8, 6).Value = 500
$ws.Cells.Item(8, 7).Value = 168
$ws.Cells.Item(8, 8).Value = 668
$ws.Cells.Item(8, 9).Value = 'cod'
$ws.Cells.Item(8, 10).Value = 'completed'
$ws.Cells.Item(8, 11).Value = '[{"product_id": 1, "product_name": "Eterno Shirt", "quantity": 1, "price": 500.0}]'
$ws.Cells.Item(8, 12).Value = '2025-11-09 12:13:27'
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 3).Value = 'kaizen'
$ws.Cells.Item(9, 4).Value = 'boarratjabol@gmail.com'
$ws.Cells.Item(9, 5).Value = 'N/A'
$ws.Cells.Item(9, 6).Value = 4599
$ws.Cells.Item(9, 7).Value = 76
$ws.Cells.Item(9, 8).Value = 4575
$ws.Cells.Item(9, 9).Value = 'credit_card'
$ws.Cells.Item(9, 10).Value = 'completed'
$ws.Cells.Item(9, 11).Value = '[{"product_id": 8, "product_name": "Eterno Abyss", "quantity": 1, "price": 4599.0}]'
$ws.Cells.Item(9, 12).Value = '2025-11-09 14:01:41'
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 3
$ws.Cells.Item(10, 3).Value = 'kaizen'
$ws.Cells.Item(10, 4).Value = 'boarratjabol@gmail.com'
$ws.Cells.Item(10, 5).Value = 'N/A'
$ws.Cells.Item(10, 6).Value = 25799
$ws.Cells.Item(10, 7).Value = 50
$ws.Cells.Item(10, 8).Value = 25749
$ws.Cells.Item(10, 9).Value = 'credit_card'
$ws.Cells.Item(10, 10).Value = 'completed'
$ws.Cells.Item(10, 11).Value = '[{"product_id": 11, "product_name": "Bonest Gatti", "quantity": 1, "price": 25799.0}]'
$ws.Cells.Item(10, 12).Value = '2025-11-09 14:16:10'

# ---- Save ----
$wb.Save()
